$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "sd_ETR"

$values = @(
    0.0228528582314404,
    69.7516175337912,
    62.9941224009146,
    14.5813819048533,
    0,
    57.5029398904697,
    57.0544687762764,
    4.16624372276806,
    0,
    73.4574235772548,
    57.7277563340043,
    0.177694826962814,
    0.274845561500558,
    53.6067191638596,
    27.603818621194,
    0
)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 6).Value = $v
    $row = $row + 1
}
